$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two new columns at D:E, shifting existing quarterly data to F:M
$ws.Columns("D:E").Insert()

# 2) Copy number formatting from column F (the shifted former column D) into the new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3) Populate the two new period columns (D = latest quarter, E = prior quarter) with reported financial figures
$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 132600
$ws.Cells.Item(8, 5).Value = 133700
$ws.Cells.Item(9, 4).Value = 109900
$ws.Cells.Item(9, 5).Value = 91500
$ws.Cells.Item(10, 4).Value = 22700
$ws.Cells.Item(10, 5).Value = 42200
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(17, 4).Value = 123600
$ws.Cells.Item(17, 5).Value = 99700
$ws.Cells.Item(18, 4).Value = 9000
$ws.Cells.Item(18, 5).Value = 34000
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(21, 4).Value = "NA"
$ws.Cells.Item(21, 5).Value = "NA"
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 4).Value = 9000
$ws.Cells.Item(23, 5).Value = 34000
$ws.Cells.Item(24, 4).Value = -500
$ws.Cells.Item(24, 5).Value = -500
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 9500
$ws.Cells.Item(26, 5).Value = 34400
$ws.Cells.Item(27, 4).Value = 3700
$ws.Cells.Item(27, 5).Value = 28000
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(33, 4).Value = 3700
$ws.Cells.Item(33, 5).Value = 28000
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 3700
$ws.Cells.Item(35, 5).Value = 28000
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 103700
$ws.Cells.Item(41, 5).Value = 57500
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(43, 4).Value = 109400
$ws.Cells.Item(43, 5).Value = 94400
$ws.Cells.Item(44, 4).Value = 29700
$ws.Cells.Item(44, 5).Value = 29600
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(47, 4).Value = 14288300
$ws.Cells.Item(47, 5).Value = 11905000
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(49, 4).Value = 25200
$ws.Cells.Item(49, 5).Value = 25200
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 5400
$ws.Cells.Item(52, 5).Value = 8300
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 14737600
$ws.Cells.Item(54, 5).Value = 12326900
$ws.Cells.Item(57, 4).Value = "NA"
$ws.Cells.Item(57, 5).Value = "NA"
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(59, 4).Value = 101200
$ws.Cells.Item(59, 5).Value = 88200
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(61, 4).Value = 13456100
$ws.Cells.Item(61, 5).Value = 11130000
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 13558200
$ws.Cells.Item(66, 5).Value = 11219200
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 289800
$ws.Cells.Item(70, 5).Value = 289800
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = -103200
$ws.Cells.Item(72, 5).Value = -75700
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 889600
$ws.Cells.Item(76, 5).Value = 817900
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 3700
$ws.Cells.Item(81, 5).Value = 28000
$ws.Cells.Item(83, 4).Value = 0
$ws.Cells.Item(83, 5).Value = 0
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 5100
$ws.Cells.Item(89, 5).Value = 8200
$ws.Cells.Item(91, 4).Value = -100
$ws.Cells.Item(91, 5).Value = -100
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -783700
$ws.Cells.Item(94, 5).Value = 8900
$ws.Cells.Item(96, 4).Value = -28200
$ws.Cells.Item(96, 5).Value = -24900
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = 822000
$ws.Cells.Item(100, 5).Value = -41700
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(102, 4).Value = 43400
$ws.Cells.Item(102, 5).Value = -24600

# 4) A handful of historical quarters were restated; correct the shifted values
$ws.Cells.Item(44, 6).Value = 29500
$ws.Cells.Item(44, 7).Value = 29300
$ws.Cells.Item(44, 8).Value = 64200
$ws.Cells.Item(44, 9).Value = 64100
$ws.Cells.Item(44, 10).Value = 34800
$ws.Cells.Item(47, 6).Value = 11214900
$ws.Cells.Item(47, 7).Value = 11383000
$ws.Cells.Item(47, 8).Value = 11614100
$ws.Cells.Item(47, 9).Value = 9643800
$ws.Cells.Item(47, 10).Value = 9797500
$ws.Cells.Item(70, 7).Value = 159300
$ws.Cells.Item(70, 9).Value = 289800
$ws.Cells.Item(91, 6).Value = 0
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 0
